# "Fix id error on diagram" — add the "Out of scope" column to the
# Software Component table (Table1) on the first worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Software Component")

# Grow Table1 by one column (auto-expands the table range + header row).
$tbl = $ws.ListObjects.Item(1)
$col = $tbl.ListColumns.Add()

# Header + data for the new column.
$ws.Cells.Item(1, 4).Value = "Out of scope"
$ws.Cells.Item(2, 4).Value = "No"
$ws.Cells.Item(3, 4).Value = "Yes"
$ws.Cells.Item(4, 4).Value = "Yes"
$ws.Cells.Item(5, 4).Value = "Yes"

# Match the look of the other bestFit columns next to it.
$ws.Columns.Item(4).AutoFit() | Out-Null
